$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Edit the hyperlink on B2: keep the same target address (mailto:grps@823)
# but change the displayed text to "grps823".
$link = $ws.Hyperlinks.Item(1)
$link.TextToDisplay = "grps823"

# Move the selection to B2 (where the edit happened)
$ws.Range("B2").Select()
